# Update cryptocurrency price/volume columns (D, E) for rows 2-51
# as scraped by the periodic GitHub Actions refresh job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "76.069.58"
$ws.Range("E2").Value = "  +1.73%  "
$ws.Range("D3").Value = "2.920.08"
$ws.Range("E3").Value = "  +2.84%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'199.59"
$ws.Range("E5").Value = "  +5.91%  "
$ws.Range("D6").Value = "'599.35"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  -1.30%  "
$ws.Range("E9").Value = "  +1.87%  "
$ws.Range("D10").Value = "2.918.81"
$ws.Range("E10").Value = "  +2.87%  "
$ws.Range("D11").Value = "'0.429"
$ws.Range("E11").Value = "  +16.05%  "
$ws.Range("E12").Value = "  -0.80%  "
$ws.Range("D13").Value = "'4.89"
$ws.Range("E13").Value = "  +0.06%  "
$ws.Range("D14").Value = "3.457.25"
$ws.Range("E14").Value = "  +2.73%  "
$ws.Range("D15").Value = "75.911.57"
$ws.Range("E15").Value = "  +1.19%  "
$ws.Range("E16").Value = "  +1.93%  "
$ws.Range("D17").Value = "'27.43"
$ws.Range("E17").Value = "  +1.06%  "
$ws.Range("D18").Value = "2.917.22"
$ws.Range("E18").Value = "  +2.30%  "
$ws.Range("E19").Value = "  -1.79%  "
$ws.Range("D20").Value = "'12.80"
$ws.Range("E20").Value = "  +3.05%  "
$ws.Range("E21").Value = "  +0.66%  "
$ws.Range("E22").Value = "  +3.01%  "
$ws.Range("E23").Value = "  +1.45%  "
$ws.Range("D24").Value = "'71.45"
$ws.Range("E24").Value = "  +0.93%  "
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("D26").Value = "3.069.94"
$ws.Range("E26").Value = "  +2.90%  "
$ws.Range("D27").Value = "'4.24"
$ws.Range("E27").Value = "  +0.83%  "
$ws.Range("D28").Value = "'9.75"
$ws.Range("E28").Value = "  +1.43%  "
$ws.Range("E29").Value = "  +5.51%  "
$ws.Range("E30").Value = "  +0.28%  "
$ws.Range("D31").Value = "'1.41"
$ws.Range("E31").Value = "  +0.43%  "
$ws.Range("D32").Value = "'506.09"
$ws.Range("E32").Value = "  -4.57%  "
$ws.Range("E33").Value = "  -2.60%  "
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").Value = "'164.90"
$ws.Range("E36").Value = "  +1.73%  "
$ws.Range("D37").Value = "'20.20"
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("E38").Value = "  +2.15%  "
$ws.Range("D39").Value = "'0.113"
$ws.Range("E39").Value = "  -5.85%  "
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("D41").Value = "'180.30"
$ws.Range("E41").Value = "  -2.04%  "
$ws.Range("E42").Value = "  +0.70%  "
$ws.Range("E43").Value = "  -1.37%  "
$ws.Range("D44").Value = "'0.0952"
$ws.Range("E44").Value = "  +10.32%  "
$ws.Range("E45").Value = "  -1.69%  "
$ws.Range("D46").Value = "'40.21"
$ws.Range("E46").Value = "  +1.46%  "
$ws.Range("D47").Value = "'1.21"
$ws.Range("E47").Value = "  -2.47%  "
$ws.Range("E48").Value = "  -0.69%  "
$ws.Range("D49").Value = "'0.578"
$ws.Range("E49").Value = "  +0.88%  "
$ws.Range("D50").Value = "'0.662"
$ws.Range("E50").Value = "  +7.37%  "
$ws.Range("E51").Value = "  -0.88%  "
